# Append the latest daily profit row (row 40) to Sheet1, matching the
# existing pattern: column A is the date stored as plain text (e.g.
# "09/26/2025"), column B is the numeric profit value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date to be written as literal text rather than having Excel's
# autocorrect turn "09/26/2025" into a date serial number: mark the cell as
# text first, assign the value, then clear the formatting so the cell keeps
# the default (unstyled) look used by the other date cells in the column.
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = "09/26/2025"
$ws.Range("A40").ClearFormats()

$ws.Range("B40").Value = 14664.18
